# Colour the "Requirements" checklist bullet items red (FF0000) to mark
# the features that are now implemented, per the commit message:
# "Implemented Dynamic product loading, search and filter functionalities
#  with Javascript".
#
# For the numbered/bulleted list paragraphs, Word stores run-level
# formatting both on the run itself (w:r/w:rPr) and on the paragraph
# mark (w:pPr/w:rPr) -- setting Font.Color on the full Paragraph.Range
# (which includes the trailing paragraph-mark character) updates both.
# For the two paragraphs where only the visible run should turn red
# (paragraph mark left untouched), a range trimmed to exclude the final
# paragraph-mark character is used instead.

$d = $word.ActiveDocument

$wdRed = 255  # RGB(255,0,0) packed as 0x0000FF -> OOXML w:color "FF0000"

# Paragraphs whose run AND paragraph mark (pPr/rPr) both get coloured red.
$fullParaTargets = @(
    "Make search bar work",
    "Results should be loaded with JavaScript",
    "Filtering features should be working",
    "Update price instantly based on the selected quantity"
)

# Paragraphs where only the run text (not the paragraph mark) gets coloured red.
$runOnlyTargets = @(
    "Sorting features should be working",
    "Loading product data with JavaScript"
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    $text = $full.Text.TrimEnd([char]13, [char]7)

    if ($fullParaTargets -contains $text) {
        $full.Font.Color = $wdRed
    }
    elseif ($runOnlyTargets -contains $text) {
        $trimmed = $d.Range($full.Start, $full.End - 1)
        $trimmed.Font.Color = $wdRed
    }
}
